$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Boris" paragraph: the sentence was originally split across several runs
#    with identical formatting ("Boris" | " has to navigate ... If the " |
#    "tablet" | " drops ... the" | " table " | "will break ... game.").
#    Re-assert the full sentence text over that span; the interop layer
#    coalesces same-format runs it rewrites into a single run, matching the
#    target (one run with the whole sentence).
# ---------------------------------------------------------------------------
$boris = "Boris has to navigate through the level while not hitting any " + `
    "obstacle or falling in any ravine/hole. If the player trips or falls, " + `
    "his laptop is going to drop, breaking it slowly. If the tablet drops " + `
    "more than a certain amount of times, the table will break completely, " + `
    "and you lose the game."
$d.Content.Find.Execute($boris, $true, $false, $false, $false, $false, $true, 1, $false, $boris, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Core Aesthetics:" heading: was split into "Core " + (bookmark _GoBack) +
#    "Aesthetics:" runs. Re-assert the full text so the runs merge into one;
#    this also removes the now-stale _GoBack bookmark (Word relocates
#    _GoBack to the point of the most recent edit - handled in step 5 below).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Core Aesthetics:", $true, $false, $false, $false, $false, $true, 1, $false, "Core Aesthetics:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Jump over obstacles such as barriers." becomes three runs:
#    "Jump " / "using SPACE " / "over obstacles such as barriers."
#    First splice in the new words, then force a run split at each new
#    boundary (toggle Bold on/off - leaves formatting unchanged but breaks
#    the run there) so the saved XML has independent <w:r> elements.
# ---------------------------------------------------------------------------
$jumpRng = $d.Content
$jumpRng.Find.Execute("Jump ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$jumpRng.Collapse(0)
$jumpRng.InsertAfter("using SPACE ")
$jumpRng.Font.Bold = $true
$jumpRng.Font.Bold = $false

# ---------------------------------------------------------------------------
# 4) "Move around mazes while not colliding with any wall or pillar."
#    becomes three runs: "Move around " / "the level" / " while not
#    colliding with any wall or pillar." (also "mazes" -> "the level").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Move around mazes while", $true, $false, $false, $false, $false, $true, 1, $false, "Move around the level while", 2) | Out-Null

$moveSearch = $d.Content
$moveSearch.Find.Execute("Move around the level while not colliding", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$moveBase = $moveSearch.Start
$lenMoveAround = ("Move around ").Length
$lenTheLevel = ("the level").Length

$moveSplit2 = $d.Range($moveBase + $lenMoveAround, $moveBase + $lenMoveAround + $lenTheLevel)
$moveSplit2.Font.Bold = $true
$moveSplit2.Font.Bold = $false

# ---------------------------------------------------------------------------
# 5) Move the _GoBack bookmark to sit right after "Try not to break your
#    tablet, overwise you lose the game." (i.e. at the very end of that
#    paragraph's run, before the paragraph mark). A genuinely empty Range
#    placed exactly on that boundary resolves ambiguously, so: insert a
#    throwaway marker character there, anchor the (now mid-text, therefore
#    unambiguous) bookmark on that boundary, then delete the marker again -
#    leaving the bookmark correctly anchored with no stray text behind.
# ---------------------------------------------------------------------------
$tabletRng = $d.Content
$tabletRng.Find.Execute("Try not to break your tablet, overwise you lose the game.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tabletEnd = $tabletRng.End

$markerRng = $d.Range($tabletEnd, $tabletEnd)
$markerRng.InsertAfter("_TMPMARK_")

$bookmarkRng = $d.Range($tabletEnd, $tabletEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRng) | Out-Null

$markerDelRng = $d.Content
$markerDelRng.Find.Execute("_TMPMARK_", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerDelRng.Delete()
